# Auto-generated Excel COM-interop script
# Applies the data updates described by the commit diff across 8 sheets
$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 7007
$ws.Range("I20").Value = 7007
$ws.Range("K20").Value = 7007
$ws.Range("M20").Value = -6777
$ws.Range("H35").Value = 7007
$ws.Range("I35").Value = 7007
$ws.Range("K35").Value = 7007
$ws.Range("M35").Value = -6628
$ws.Range("H137").Value = 2299.9
$ws.Range("I137").Value = 1570.8572
$ws.Range("J137").Value = 4001
$ws.Range("K137").Value = 4712.571599999999
$ws.Range("L137").Value = 12003
$ws.Range("M137").Value = -2162.571599999999
$ws.Range("N137").Value = -17103
$ws.Range("H141").Value = 3130
$ws.Range("I141").Value = 1494.3103
$ws.Range("J141").Value = 8400.556
$ws.Range("K141").Value = 4482.9309
$ws.Range("L141").Value = 25201.668
$ws.Range("M141").Value = 697.0690999999997
$ws.Range("N141").Value = -35561.66800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H35").Value = 1500
$ws.Range("I35").Value = 1500
$ws.Range("K35").Value = 1500
$ws.Range("M35").Value = -1094
$ws.Range("H39").Value = 8033.3335
$ws.Range("I39").Value = 8033.3335
$ws.Range("K39").Value = 8033.3335
$ws.Range("M39").Value = -7513.3335
$ws.Range("H45").Value = 60146.94
$ws.Range("I45").Value = 91969.17999999999
$ws.Range("J45").Value = 1806.1666
$ws.Range("K45").Value = 91969.17999999999
$ws.Range("L45").Value = 1806.1666
$ws.Range("M45").Value = -91592.17999999999
$ws.Range("N45").Value = -2560.1666
$ws.Range("H110").Value = 1597.1904
$ws.Range("I110").Value = 820.05884
$ws.Range("J110").Value = 4900
$ws.Range("K110").Value = 820.05884
$ws.Range("L110").Value = 4900
$ws.Range("M110").Value = 1224.94116
$ws.Range("N110").Value = -8990
$ws.Range("H132").Value = 1032974.5
$ws.Range("I132").Value = 786.5833
$ws.Range("J132").Value = 6537976.5
$ws.Range("K132").Value = 2359.7499
$ws.Range("L132").Value = 19613929.5
$ws.Range("M132").Value = 170.2501000000002
$ws.Range("N132").Value = -19618989.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 799.4
$ws.Range("I94").Value = 684.1429000000001
$ws.Range("J94").Value = 900.25
$ws.Range("K94").Value = 684.1429000000001
$ws.Range("L94").Value = 900.25
$ws.Range("M94").Value = -233.1429000000001
$ws.Range("N94").Value = -1802.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 902.5
$ws.Range("I12").Value = 902.5
$ws.Range("K12").Value = 902.5
$ws.Range("M12").Value = -732.5
$ws.Range("H25").Value = 2655.5
$ws.Range("I25").Value = 2655.5
$ws.Range("K25").Value = 2655.5
$ws.Range("M25").Value = -2481.5
$ws.Range("H31").Value = 1404.48
$ws.Range("I31").Value = 913.26666
$ws.Range("J31").Value = 2141.3
$ws.Range("K31").Value = 913.26666
$ws.Range("L31").Value = 2141.3
$ws.Range("M31").Value = -618.26666
$ws.Range("N31").Value = -2731.3
$ws.Range("H34").Value = 1404.48
$ws.Range("I34").Value = 913.26666
$ws.Range("J34").Value = 2141.3
$ws.Range("K34").Value = 913.26666
$ws.Range("L34").Value = 2141.3
$ws.Range("M34").Value = -711.26666
$ws.Range("N34").Value = -2545.3
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H58").Value = 21739708
$ws.Range("I58").Value = 30303492
$ws.Range("J58").Value = 873.0769
$ws.Range("K58").Value = 30303492
$ws.Range("L58").Value = 873.0769
$ws.Range("M58").Value = -30303289
$ws.Range("N58").Value = -1279.0769
$ws.Range("H132").Value = 10102663
$ws.Range("I132").Value = 1170.6842
$ws.Range("J132").Value = 23811832
$ws.Range("K132").Value = 3512.0526
$ws.Range("L132").Value = 71435496
$ws.Range("M132").Value = -982.0526
$ws.Range("N132").Value = -71440556
$ws.Range("H136").Value = 21739708
$ws.Range("I136").Value = 30303492
$ws.Range("J136").Value = 873.0769
$ws.Range("K136").Value = 90910476
$ws.Range("L136").Value = 2619.2307
$ws.Range("M136").Value = -90907926
$ws.Range("N136").Value = -7719.2307

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 250362.75
$ws.Range("I4").Value = 333717
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1001151
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = -1001039
$ws.Range("N4").Value = -1124
$ws.Range("H17").Value = 322.5
$ws.Range("I17").Value = 96.666664
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 289.999992
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -120.999992
$ws.Range("N17").Value = -3338
$ws.Range("H38").Value = 75.111115
$ws.Range("I38").Value = 40.5
$ws.Range("J38").Value = 102.8
$ws.Range("K38").Value = 121.5
$ws.Range("L38").Value = 308.4
$ws.Range("M38").Value = 225.5
$ws.Range("N38").Value = -1002.4
$ws.Range("H131").Value = 908.78
$ws.Range("J131").Value = 925.5
$ws.Range("L131").Value = 2776.5
$ws.Range("N131").Value = -12856.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 800
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -384
$ws.Range("N3").Value = -1232
$ws.Range("H7").Value = 3002000
$ws.Range("I7").Value = 5000000
$ws.Range("J7").Value = 1670000
$ws.Range("K7").Value = 5000000
$ws.Range("L7").Value = 1670000
$ws.Range("M7").Value = -4999888
$ws.Range("N7").Value = -1670224
$ws.Range("H8").Value = 3002000
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 1670000
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 1670000
$ws.Range("M8").Value = -4999861
$ws.Range("N8").Value = -1670278
$ws.Range("H12").Value = 2627870
$ws.Range("H14").Value = 172.375
$ws.Range("I14").Value = 172.375
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 172.375
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -4.375
$ws.Range("N14").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 2750.5
$ws.Range("I12").Value = 501
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 501
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -331
$ws.Range("N12").Value = -5340
$ws.Range("H30").Value = 2258.8
$ws.Range("I30").Value = 758.6667
$ws.Range("J30").Value = 4509
$ws.Range("K30").Value = 758.6667
$ws.Range("L30").Value = 4509
$ws.Range("M30").Value = -650.6667
$ws.Range("N30").Value = -4725
$ws.Range("H93").Value = 1087.5
$ws.Range("J93").Value = 1144.5555
$ws.Range("L93").Value = 1144.5555
$ws.Range("N93").Value = -3640.5555
$ws.Range("H101").Value = 27362
$ws.Range("J101").Value = 27362
$ws.Range("L101").Value = 27362
$ws.Range("N101").Value = -33852

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H28").Value = 1908.5
$ws.Range("J28").Value = 2800
$ws.Range("L28").Value = 2800
$ws.Range("N28").Value = -3496
$ws.Range("H40").Value = 3700
$ws.Range("I40").Value = 3700
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3700
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3551
$ws.Range("N40").ClearContents()
$ws.Range("H132").Value = 21532.105
$ws.Range("I132").Value = 28661.658
$ws.Range("J132").Value = 7273
$ws.Range("K132").Value = 85984.974
$ws.Range("L132").Value = 21819
$ws.Range("M132").Value = -83454.974
$ws.Range("N132").Value = -26879
